$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 3) to the meta-sheet with the PF/1.0.6 entry
$ws.Range("A3").Value = "PF/1.0.6"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# The new row uses the default/normal style (not the header style inherited
# from the column formatting), matching the rest of the sheet's data rows.
$ws.Range("A3:D3").Style = "Normal"
